$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C values per row (2..15), added as part of "Pre-tx phase self report outcome measures"
$values = @{
    2  = "Somewhat worse"
    3  = "Somewhat worse"
    4  = "A lot worse"
    5  = "A lot worse"
    6  = "A lot worse"
    7  = "A lot worse"
    8  = "Somewhat worse"
    9  = "A little worse"
    10 = "A little worse"
    11 = "Somewhat worse"
    12 = "Somewhat worse"
    13 = "A lot worse"
    14 = "A lot worse"
    15 = "A lot worse"
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 3).Value = $values[$row]
}

# Update the active selection to reflect where the user ended up editing
$ws.Range("C16").Select()
